# Update sheet "voice_live_environment": replace "mean_amplitude"/"max_local_amplitude"/"file_name"
# 3-column report (time_(s), mean_amplitude, max_local_amplitude, file_name; 21 data rows) with a
# new 2-measurement report (time_(s), max_local_amplitude, file_name; 18 data rows) covering the
# "environment_noise" recordings only ("complete system no noise").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("voice_live_environment")

# Clear the previous table's values first (A2:D22), keeping header formatting intact for A2:C2.
$ws.Range("A2:D22").ClearContents()
# Column D is being dropped entirely, so wipe its formatting too.
$ws.Range("D2:D22").Clear()

# Header row.
$ws.Cells.Item(2, 1).Value = "time_(s)"
$ws.Cells.Item(2, 2).Value = "max_local_amplitude"
$ws.Cells.Item(2, 3).Value = "file_name"

$data = @(
    @(1.954366445541382, 0.04589070007205009, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191413455.wav'),
    @(2.959033727645874, 0.05773304402828217, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191413460.wav'),
    @(3.962219476699829, 0.04831305518746376, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914134713.wav'),
    @(4.9663405418396, 0.0461459681391716, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914134816.wav'),
    @(5.969293832778931, 0.04763097316026688, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914134914.wav'),
    @(6.971488237380981, 0.1183869242668152, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191413509.wav'),
    @(7.97303318977356, 0.1052666828036308, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191413518.wav'),
    @(8.974731683731079, 0.0661206841468811, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191413529.wav'),
    @(9.977128982543945, 0.05795527622103691, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191413530.wav'),
    @(10.97979760169983, 0.05269449204206467, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914135410.wav'),
    @(11.98143911361694, 0.09235147386789322, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914135514.wav'),
    @(12.98377537727356, 0.07962875068187714, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914135610.wav'),
    @(13.98786497116089, 0.05643543228507042, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914135719.wav'),
    @(14.9930534362793, 0.0741506814956665, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914135810.wav'),
    @(15.99627375602722, 0.0663648471236229, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/2021051914135917.wav'),
    @(17.00104832649231, 0.09688328206539154, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191414004.wav'),
    @(18.00464606285095, 0.08851940184831619, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191414017.wav'),
    @(19.00777387619019, 0.127306342124939, '/home/minhair/Desktop/food_ordering_system/food_ordering_system/predict/recorded_audios/environment_noise/202105191414025.wav'),
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

